$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "ShortName"
$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("C5").Select()
